$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean so shared strings / cells get rebuilt from scratch.
$ws.Cells.Clear()

# Re-enter the cell values in the same order they were authored so the
# resulting shared-strings table lines up with the target workbook.
$ws.Range("B9").Value  = "Der User kann sein Konto löschen. (Sicherheit)"
$ws.Range("B18").Value = "jQuery"
$ws.Range("B17").Value = "JavaScript"
$ws.Range("B20").Value = "php"
$ws.Range("A16").Value = "Verwendete Technologien:"
$ws.Range("B21").Value = "(SpotifyAPI)"
$ws.Range("A3").Value  = "Antwortzeiten"
$ws.Range("A8").Value  = "Sicherheit"
$ws.Range("A12").Value = "Usability"
$ws.Range("B13").Value = "Das Game benötigt keine Anleitung."
$ws.Range("A1").Value  = "Nicht funktionale Anforderungen an SongQuiz"
$ws.Range("B4").Value  = "Beim Beantworten dauert die Anzeige ob richtig oder falsch maximal eine Sekunde."
$ws.Range("B5").Value  = "Das Starten eines neuen Spiels dauert nicht länger als 5 Sekunden (wird sich noch zeigen je nach API Funktionalität)"
$ws.Range("B19").Value = "MySQL"

# Column widths: column A narrow, column B wide enough to fit the long
# requirement text (mirrors the bestFit/customWidth columns in the target).
$ws.Columns.Item(1).ColumnWidth = 2.75
$ws.Columns.Item(2).ColumnWidth = 93.75

# Selection moves to B6 in the saved view.
$null = $ws.Range("B6").Select()
